# CAFÉ ARENILLO.docx - schedule table update
# - 08:30-08:45 slot buyer changes from FLOR A FRUTO to INTERLINK2AMERICAS
# - New row inserted right after it: 08:45-09:00 / FLOR A FRUTO
# - COLFRESH COFFEE slot time 09:45-10:00 -> 09:00-09:15
# - Old INTERLINK2AMERICAS row (10:00-10:15) removed (it moved to the first slot)
# - INMERSSO BOUTIQUE slot time 10:30-10:45 -> 10:00-10:15
# - BOX BRAND slot time 11:30-11:45 -> 10:15-10:30

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Step 1: Row 2 (08:30 - 08:45) buyer FLOR A FRUTO -> INTERLINK2AMERICAS
$t.Rows.Item(2).Cells.Item(3).Range.Text = "INTERLINK2AMERICAS"

# Step 2: Insert a new row right after row 2 (before the current row 3,
# which is 09:45-10:00 / COLFRESH COFFEE). Populate with 08:45-09:00 / FLOR A FRUTO
$newRow = $t.Rows.Add($t.Rows.Item(3))
$newRow.Cells.Item(1).Range.Text = "08:45 - 09:00"
$newRow.Cells.Item(3).Range.Text = "FLOR A FRUTO"

# Step 3: Row 4 is now COLFRESH COFFEE (09:45 - 10:00) -> 09:00 - 09:15
$t.Rows.Item(4).Cells.Item(1).Range.Text = "09:00 - 09:15"

# Step 4: Row 5 is now the old INTERLINK2AMERICAS row (10:00 - 10:15) -> delete it
$t.Rows.Item(5).Delete()

# Step 5: Row 5 is now INMERSSO BOUTIQUE (10:30 - 10:45) -> 10:00 - 10:15
$t.Rows.Item(5).Cells.Item(1).Range.Text = "10:00 - 10:15"

# Step 6: Row 6 is now BOX BRAND (11:30 - 11:45) -> 10:15 - 10:30
$t.Rows.Item(6).Cells.Item(1).Range.Text = "10:15 - 10:30"
